$wb = $excel.ActiveWorkbook
$sheet1 = $wb.Worksheets.Item(1)

# --- Add the new "Classlist" worksheet after Sheet1 ---
$classlist = $wb.Worksheets.Add($null, $sheet1)
$classlist.Name = "Classlist"

# --- Populate the class list table (A3:D8) ---
# Row 3 filled completely first.
$classlist.Range("A3").Value = "Mbongwe,"
$classlist.Range("B3").Value = "KA,"
$classlist.Range("C3").Value = "Miss"
$classlist.Range("D3").Value = "[22639438@sun.ac.za]"

# Row 4: B,C,D filled now; A4 ("Du Toit,") is filled in later (matches
# the original author's edit order / shared-string insertion order).
$classlist.Range("B4").Value = "F,"
$classlist.Range("C4").Value = "Mnr"
$classlist.Range("D4").Value = "[22552987@sun.ac.za]"

$classlist.Range("A5").Value = "Matthysen,"
$classlist.Range("B5").Value = "LP,"
$classlist.Range("C5").Value = "Mr"
$classlist.Range("D5").Value = "[22899537@sun.ac.za]"

$classlist.Range("A6").Value = "Mofokeng,"
$classlist.Range("B6").Value = "D,"
$classlist.Range("C6").Value = "Me"
$classlist.Range("D6").Value = "[22309667@sun.ac.za]"

$classlist.Range("A7").Value = "Schultz,"
$classlist.Range("B7").Value = "K,"
$classlist.Range("C7").Value = "Mnr"
$classlist.Range("D7").Value = "[22539026@sun.ac.za]"

$classlist.Range("A8").Value = "Sinclair,"
$classlist.Range("B8").Value = "HW,"
$classlist.Range("C8").Value = "Mr"
$classlist.Range("D8").Value = "[21672598@sun.ac.za]"

# A4 written last.
$classlist.Range("A4").Value = "Du Toit,"

# --- Column widths on the Classlist sheet ---
$classlist.Columns.Item(1).ColumnWidth = 36.498697916666664
$classlist.Columns.Item(4).ColumnWidth = 19.166666666666668

# --- View / selection state ---
# Sheet1: zoom to 170%, selection moves to J9, no longer the tab in focus.
$sheet1.Select()
$excel.ActiveWindow.Zoom = 170
$sheet1.Range("J9").Select() | Out-Null

# Classlist: becomes the active/selected sheet, selection at E20.
$classlist.Select()
$classlist.Range("E20").Select() | Out-Null
